$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.333.73"
$ws.Range("D3").Value = "1.832.66"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.90%  "
$ws.Range("D5").Value = "'314.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").Value = "'0.4749"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("D8").Value = "'0.3686"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "'0.07458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'0.8852"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "1.883.23"
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("D13").Value = "'0.07336"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "'5.441"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'93.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "'6.579"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.000008796"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "27.552.15"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "'5.292"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D24").Value = "2.096.00"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("D25").Value = "'1.890"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'151.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").Value = "'18.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "'5.247"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'117.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "'0.7525"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "'2.944"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "'1.011"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").Value = "'0.05347"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'2.978"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "'2.391"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("D43").Value = "'0.5314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D45").Value = "'8.476"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "'0.4913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'10.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D49").Value = "'105.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'0.06299"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.16%  "
